# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: header cell "Save", styled like the other header cells (copy format from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# H2 / H3: numeric data cells, value 0 (no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
